$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
$ws.Columns(8).Delete()
$wb.Styles("Heading 2").Delete()
